$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new cell F3 = "someValue"
$ws.Range("F3").Value = "someValue"

# Add new row 7 with A7 = "test for row"
$ws.Range("A7").Value = "test for row"

# Update the selection to match the post-edit active cell (A8)
$ws.Range("A8").Select()
